# Apply edits for "Finish Hotels Search test case & Enhancements"
$wb = $excel.ActiveWorkbook

# The relevant worksheet (currently named "testsheet2") is the active sheet.
$ws = $wb.ActiveSheet

# Rename it to "GUI"
$ws.Name = "GUI"

# Add the new "Expected Hotel Name" column header and its expected value
$ws.Range("E1").Value = "Expected Hotel Name"
$ws.Range("E2").Value = "Grand Plaza Apartments"

# Copy header style from D1 to E1, and data style from D2 to E2
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)  # xlPasteFormats

# Widen the new column to fit its (longer) content, mirroring the other
# "best fit" columns on this sheet
$ws.Columns.Item(5).ColumnWidth = 20.3

# Update the active selection to E2
$ws.Range("E2").Select()
